# Ajusted mvp-requisitos, cadastre-se form and fixed avaliacao.html button
#
# The original "requisitos" table had two requirement rows that were
# duplicates / out of place (row 14 "pagina de usuario com login" text and
# row 10 "qualquer pessoa posso se cadastrar" typo) as well as a whole
# requirement row (old row 11, "relatorios em pdf") and another one (old
# row 18 duplicate) that got merged/removed. The sheet below is rewritten
# so that the remaining 16 requirement rows (2-17) hold the corrected /
# reordered text, the row heights for the wrapped cells match the new
# text length, and the now-unused trailing row 18 is deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : Não Funcional / disponibilidade -------------------------------
$ws.Range("A2").Value = "Não Funcional"
$ws.Range("B2").Value = "O sistema deve ficar disponível, no mínimo, 95% do tempo anual, eventuais manutenções deverão ocorrer de preferência no período noturno de menor atividade;"
$ws.Range("C2").Value = "Importante"

# --- Row 3 : Não Funcional / tempo de resposta ------------------------------
$ws.Range("A3").Value = "Não Funcional"
$ws.Range("B3").Value = "O sistema deve ter um tempo de resposta inferior a 4 segundos;"
$ws.Range("C3").Value = "Importante"

# --- Row 4 : Não Funcional / navegadores ------------------------------------
$ws.Range("A4").Value = "Não Funcional"
$ws.Range("B4").Value = "O sistema deve ser compatível com os navegadores Chrome, Edge, Opera e Firefox;"
$ws.Range("C4").Value = "Desejável"

# --- Row 5 : Não Funcional / responsivo -------------------------------------
$ws.Range("A5").Value = "Não Funcional"
$ws.Range("B5").Value = "O sistema deve ser responsivo para smartphones;"
$ws.Range("C5").Value = "Essencial"

# --- Row 6 : Não Funcional / intuitivo --------------------------------------
$ws.Range("A6").Value = "Não Funcional"
$ws.Range("B6").Value = "O sistema deve ser ituitivo e fácil de usar, com cores confortáveis e informações objetivas;"
$ws.Range("C6").Value = "Essencial"

# --- Row 7 : Não Funcional / tecnologias front-end --------------------------
$ws.Range("A7").Value = "Não Funcional"
$ws.Range("B7").Value = "O  sistema deve ser desenvolvido com as seguintes tecnologias Front-End: Html, Css e JavaScript;"
$ws.Range("C7").Value = "Essencial"

# --- Row 8 : Não Funcional / login simultâneo -------------------------------
$ws.Range("A8").Value = "Não Funcional"
$ws.Range("B8").Value = "O sistema não deve permitir que o usuário se conecte simultaneamente em aparelhos distintos  com o mesmo login;"
$ws.Range("C8").Value = "Importante"

# --- Row 9 : Não Funcional / LGPD -------------------------------------------
$ws.Range("A9").Value = "Não Funcional"
$ws.Range("B9").Value = "O sistema deve armazenar e proteger dados de acordo com a Lei Geral de Proteção de Dados."
$ws.Range("C9").Value = "Essencial"

# --- Row 13 content is written first so the two brand-new requirement
# strings land in the shared-strings table in the same order as the
# canonical file (opções de login, then qualquer pessoa possa se
# cadastrar): -----------------------------------------------------------

# --- Row 13 : Funcional / opções de login (adds "Login do cliente") --------
$ws.Range("A13").Value = "Funcional"
$ws.Range("B13").Value = "O sistema deve exibir a página de usuário com as seguintes opções de login: Login do cliente, Estudante, Porfessor e Administrador;"
$ws.Range("C13").Value = "Essencial"
$ws.Rows.Item(13).RowHeight = 88.8

# --- Row 10 : Funcional / cadastro (typo fixed: "posso" -> "possa") --------
$ws.Range("A10").Value = "Funcional"
$ws.Range("B10").Value = "O sistema deve permitir que qualquer pessoa possa se cadastrar;"
$ws.Range("C10").Value = "Essencial"

# --- Row 11 : Funcional / página principal pública --------------------------
$ws.Range("A11").Value = "Funcional"
$ws.Range("B11").Value = "O sistema deve exibir uma página principal pública com menu de navegação para endereços importantes; "
$ws.Range("C11").Value = "Desejável"

# --- Row 12 : Funcional / agendamentos de datas e horários ------------------
$ws.Range("A12").Value = "Funcional"
$ws.Range("B12").Value = "O sistema deve permitir que todas as pessoas realizem agendamentos de datas e horários disponíveis;"
$ws.Range("C12").Value = "Essencial"
$ws.Rows.Item(12).RowHeight = 102
# D12 goes from the "plain" placeholder style to the "underline" placeholder
# style used by D13/D14 - copy that formatting across (keeps the same
# shared style index instead of allocating a brand-new one).
$ws.Range("D13").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 14 : Funcional / Administrador acessa tabela -----------------------
$ws.Range("A14").Value = "Funcional"
$ws.Range("B14").Value = "O sistema deve permitir que o Administrador, com login efetuado, consiga acessar todas as informações de agendamentos em uma tabela;"
$ws.Range("C14").Value = "Essencial"
$ws.Rows.Item(14).RowHeight = 76.8
# D14 no longer has a placeholder cell at all.
$ws.Range("D14").Clear()

# --- Row 15 : Funcional / Administrador edita/deleta atendimento -----------
$ws.Range("A15").Value = "Funcional"
$ws.Range("B15").Value = "O sistema deve permitir que o Administrador, com login efeutaudo, altere/edite o grupo responsável pelo atendimento e delete um atendimento, caso seja necessário;"
$ws.Range("C15").Value = "Essencial"

# --- Row 16 : Funcional / Professores/Estudantes tabela ---------------------
$ws.Range("A16").Value = "Funcional"
$ws.Range("B16").Value = "O sistema deve permitir que Professores e Estudantes, com login efetuado, consigam acessar uma tabela com os dados dos agendamentos atribuídos à eles pelo administrador;"
$ws.Range("C16").Value = "Importante"

# --- Row 17 : Funcional / avaliação do atendimento --------------------------
$ws.Range("A17").Value = "Funcional"
$ws.Range("B17").Value = "O sistema deve permitir que um usuário logado avalie o atendimento;"
$ws.Range("C17").Value = "Desejável"
$ws.Rows.Item(17).RowHeight = 87.6

# --- Old row 18 is no longer needed - remove it, shrinking the table -------
$ws.Rows.Item(18).Delete()

# Selection ends up on B8, matching the saved view state of the edited file.
$ws.Range("B8").Select()
